$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a number by Excel;
# force them to remain text, matching the original inline-string cell type,
# then drop the temporary Text number-format so no stray style is left behind.
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D29", "D30", "D31", "D34", "D35", "D36", "D38", "D39", "D40", "D43", "D44", "D45", "D46", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.381.85'
$ws.Range("E2").Value = '  +2.49%  '
$ws.Range("D3").Value = '3.242.19'
$ws.Range("E3").Value = '  +5.04%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").Value = '577.10'
$ws.Range("E5").Value = '  +2.03%  '
$ws.Range("D6").Value = '155.33'
$ws.Range("E6").Value = '  +8.77%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = '3.231.37'
$ws.Range("E8").Value = '  +5.05%  '
$ws.Range("D9").Value = '0.516'
$ws.Range("E9").Value = '  +4.36%  '
$ws.Range("D10").Value = '7.09'
$ws.Range("E10").Value = '  +9.87%  '
$ws.Range("D11").Value = '0.167'
$ws.Range("E11").Value = '  +5.22%  '
$ws.Range("D12").Value = '0.490'
$ws.Range("E12").Value = '  +4.39%  '
$ws.Range("D13").Value = '38.16'
$ws.Range("E13").Value = '  +6.39%  '
$ws.Range("D14").Value = '0.0000236'
$ws.Range("E14").Value = '  +3.72%  '
$ws.Range("D15").Value = '3.753.60'
$ws.Range("E15").Value = '  +4.89%  '
$ws.Range("D16").Value = '557.85'
$ws.Range("E16").Value = '  +12.40%  '
$ws.Range("D17").Value = '66.382.79'
$ws.Range("E17").Value = '  +2.55%  '
$ws.Range("E18").Value = '  +3.22%  '
$ws.Range("D19").Value = '3.236.46'
$ws.Range("E19").Value = '  +4.77%  '
$ws.Range("D20").Value = '7.14'
$ws.Range("E20").Value = '  +6.36%  '
$ws.Range("D21").Value = '14.46'
$ws.Range("E21").Value = '  +4.54%  '
$ws.Range("D22").Value = '0.745'
$ws.Range("E22").Value = '  +7.45%  '
$ws.Range("D23").Value = '7.88'
$ws.Range("E23").Value = '  +8.77%  '
$ws.Range("D24").Value = '13.66'
$ws.Range("E24").Value = '  +6.87%  '
$ws.Range("D25").Value = '82.25'
$ws.Range("E25").Value = '  +4.21%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  +17.76%  '
$ws.Range("D28").Value = '2.93'
$ws.Range("E28").Value = '  +5.21%  '
$ws.Range("D29").Value = '2.28'
$ws.Range("E29").Value = '  +8.41%  '
$ws.Range("D30").Value = '27.89'
$ws.Range("E30").Value = '  +5.12%  '
$ws.Range("D31").Value = '2.78'
$ws.Range("E31").Value = '  +3.19%  '
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("E33").Value = '  +5.25%  '
$ws.Range("D34").Value = '570.06'
$ws.Range("E34").Value = '  +9.86%  '
$ws.Range("D35").Value = '5.82'
$ws.Range("E35").Value = '  +4.25%  '
$ws.Range("D36").Value = '6.45'
$ws.Range("E36").Value = '  +7.22%  '
$ws.Range("E37").Value = '  +13.34%  '
$ws.Range("D38").Value = '55.48'
$ws.Range("E38").Value = '  +3.66%  '
$ws.Range("D39").Value = '0.0874'
$ws.Range("E39").Value = '  +8.88%  '
$ws.Range("D40").Value = '3.06'
$ws.Range("E40").Value = '  +13.06%  '
$ws.Range("E41").Value = '  +4.99%  '
$ws.Range("D42").Value = '3.136.99'
$ws.Range("E42").Value = '  +6.47%  '
$ws.Range("D43").Value = '8.67'
$ws.Range("E43").Value = '  +2.97%  '
$ws.Range("D44").Value = '0.275'
$ws.Range("E44").Value = '  +10.96%  '
$ws.Range("D45").Value = '2.35'
$ws.Range("E45").Value = '  +7.85%  '
$ws.Range("D46").Value = '27.10'
$ws.Range("E46").Value = '  +6.95%  '
$ws.Range("D47").Value = '0.0₃0565'
$ws.Range("E47").Value = '  +3.39%  '
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("D49").Value = '0.114'
$ws.Range("E49").Value = '  +4.32%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '2.26'
$ws.Range("E50").Value = '  +8.76%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '122.56'
$ws.Range("E51").Value = '  +0.90%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}

